$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "caso" -> "case"
$ws.Range("A1").Value = "case"

# "fecha" -> "date"
$ws.Range("B1").Value = "date"

# B2 holds a date-like text string ("02/15/2025" -> "02/15/2026") in a cell
# formatted with a custom date numFmt ("d\-m\-yyyy;@"). Assigning a
# date-looking string directly would make Excel auto-convert it into a real
# date serial number, so temporarily switch the format to plain text,
# assign the new text, then restore the exact original custom format code.
$origFormat = "d\-m\-yyyy;@"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "02/15/2026"
$ws.Range("B2").NumberFormat = $origFormat

# Move the sheet view's active selection from B2 to D6
$ws.Range("D6").Select()
